# Scheduled runner refresh: re-pull current market averages (columns
# H/I/J/K/L) and recompute the dependent profit columns (M/N) for the
# Leve rows whose underlying item prices moved since the last run.
# Some rows lose their HQ/NQ profit cell entirely when the HQ branch no
# longer applies (and vice versa for rows that gain one).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H43").Value = 5664.4443
$ws.Range("I43").Value = 1875
$ws.Range("K43").Value = 1875
$ws.Range("M43").Value = -1806

$ws.Range("H53").Value = 87.888885
$ws.Range("I53").Value = 87.888885
$ws.Range("K53").Value = 87.888885
$ws.Range("M53").Value = 549.111115

$ws.Range("H58").Value = 536.5
$ws.Range("I58").Value = 50.333332
$ws.Range("J58").Value = 1995
$ws.Range("K58").Value = 150.999996
$ws.Range("L58").Value = 5985
$ws.Range("M58").Value = -0.9999960000000101
$ws.Range("N58").Value = -6285

$ws.Range("H74").Value = 8000
$ws.Range("I74").Value = 8000
$ws.Range("K74").Value = 8000
$ws.Range("M74").Value = -7064

$ws.Range("H77").Value = 8000
$ws.Range("I77").Value = 8000
$ws.Range("K77").Value = 40000
$ws.Range("M77").Value = -35320

$ws.Range("H99").Value = 325.875
$ws.Range("I99").Value = 223.6
$ws.Range("K99").Value = 670.8
$ws.Range("M99").Value = 827.2

$ws.Range("H107").Value = 1584.375
$ws.Range("I107").Value = 1584.375
$ws.Range("K107").Value = 1584.375
$ws.Range("M107").Value = 335.625

$ws.Range("H141").Value = 3666.3333
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 42000
$ws.Range("J54").Value = 42000
$ws.Range("L54").Value = 42000
$ws.Range("N54").Value = -43538

$ws.Range("H61").Value = 6150
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9788

$ws.Range("H102").Value = 970
$ws.Range("I102").Value = 970
$ws.Range("K102").Value = 970
$ws.Range("M102").Value = 652

$ws.Range("H136").Value = 6150
$ws.Range("I136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("M136").Value = -27450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 304.6
$ws.Range("I5").Value = 304.6
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 304.6
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -191.6
$ws.Range("N5").ClearContents()

$ws.Range("H7").Value = 980.6
$ws.Range("I7").Value = 350.75
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 350.75
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -237.75
$ws.Range("N7").Value = -3726

$ws.Range("H86").Value = 10007
$ws.Range("J86").Value = 10007
$ws.Range("L86").Value = 10007
$ws.Range("N86").Value = -12253

$ws.Range("H89").Value = 10007
$ws.Range("J89").Value = 10007
$ws.Range("L89").Value = 50035
$ws.Range("N89").Value = -61267

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 632
$ws.Range("I2").Value = 323.85715
$ws.Range("J2").Value = 1171.25
$ws.Range("K2").Value = 323.85715
$ws.Range("L2").Value = 1171.25
$ws.Range("M2").Value = -210.85715
$ws.Range("N2").Value = -1397.25

$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H6").Value = 5550745.5
$ws.Range("I6").Value = 5550745.5
$ws.Range("K6").Value = 5550745.5
$ws.Range("M6").Value = -5550632.5

$ws.Range("H7").Value = 509.42856
$ws.Range("I7").Value = 824.75
$ws.Range("K7").Value = 824.75
$ws.Range("M7").Value = -711.75

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H50").Value = 23333
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H60").Value = 7342.7144
$ws.Range("I60").Value = 7342.7144
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 7342.7144
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -6831.7144
$ws.Range("N60").ClearContents()

$ws.Range("H104").Value = 70000
$ws.Range("J104").Value = 70000
$ws.Range("L104").Value = 70000
$ws.Range("N104").Value = -75242

$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3174.6
$ws.Range("J34").Value = 4471.2856
$ws.Range("L34").Value = 13413.8568
$ws.Range("N34").Value = -13581.8568

$ws.Range("H117").Value = 3952.2
$ws.Range("I117").Value = 430
$ws.Range("K117").Value = 1290
$ws.Range("M117").Value = 2152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 1474
$ws.Range("J17").Value = 1474
$ws.Range("L17").Value = 1474
$ws.Range("N17").Value = -1810

$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 1000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -730
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 1000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -64
$ws.Range("N73").ClearContents()

$ws.Range("H97").Value = 1948.0588
$ws.Range("I97").Value = 1885.1666
$ws.Range("K97").Value = 1885.1666
$ws.Range("M97").Value = -1389.1666

$ws.Range("H102").Value = 682.8125
$ws.Range("I102").Value = 593.6429000000001
$ws.Range("J102").Value = 1307
$ws.Range("K102").Value = 593.6429000000001
$ws.Range("L102").Value = 1307
$ws.Range("M102").Value = 1028.3571
$ws.Range("N102").Value = -4551

$ws.Range("H122").Value = 1251.5555
$ws.Range("I122").Value = 1251.75
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 3755.25
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -1305.25
$ws.Range("N122").Value = -8650

$ws.Range("H126").Value = 1758.2858
$ws.Range("I126").Value = 1802.25
$ws.Range("K126").Value = 5406.75
$ws.Range("M126").Value = -2936.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -888

$ws.Range("H106").Value = 11421
$ws.Range("J106").Value = 11421
$ws.Range("L106").Value = 11421
$ws.Range("N106").Value = -13945

$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

$ws.Range("H132").Value = 2935.6
$ws.Range("I132").Value = 2935.6
$ws.Range("K132").Value = 8806.799999999999
$ws.Range("M132").Value = -6276.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
